# Weekly update: insert 3 new price rows (Especial / Primera / Segunda) for the
# new reporting date 44673 (2022-04-22) just above the previously most-recent
# entries for this market/product/variety, pushing all later rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 129 (old rows 129..203 shift down to 132..206)
$ws.Rows("129:131").Insert()

# Shared values for the new rows (same market/product/category/variety as the rest of the sheet)
$mercadoId = 4
$mercado = "Feria Lagunitas de Puerto Montt"
$region = "Los Lagos"
$codreg = 10
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100101007
$categoria = "Kiwi"
$variedad = "Hayward"
$unidad = "`$/caja 15 kilos"
$kgUnidad = 15
$fecha = 44673

# Row 129: Especial
$r = 129
$ws.Cells.Item($r, 1).Value2 = $mercadoId
$ws.Cells.Item($r, 2).Value2 = $mercado
$ws.Cells.Item($r, 3).Value2 = $region
$ws.Cells.Item($r, 4).Value2 = $fecha
$ws.Cells.Item($r, 5).Value2 = $codreg
$ws.Cells.Item($r, 6).Value2 = $tipo
$ws.Cells.Item($r, 7).Value2 = $productoId
$ws.Cells.Item($r, 8).Value2 = $producto
$ws.Cells.Item($r, 9).Value2 = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = $variedad
$ws.Cells.Item($r, 12).Value2 = "Especial"
$ws.Cells.Item($r, 13).Value2 = 300
$ws.Cells.Item($r, 14).Value2 = 19000
$ws.Cells.Item($r, 15).Value2 = 19000
$ws.Cells.Item($r, 16).Value2 = 19000
$ws.Cells.Item($r, 17).Value2 = $unidad
$ws.Cells.Item($r, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value2 = 1267
$ws.Cells.Item($r, 20).Value2 = $kgUnidad

# Row 130: Primera
$r = 130
$ws.Cells.Item($r, 1).Value2 = $mercadoId
$ws.Cells.Item($r, 2).Value2 = $mercado
$ws.Cells.Item($r, 3).Value2 = $region
$ws.Cells.Item($r, 4).Value2 = $fecha
$ws.Cells.Item($r, 5).Value2 = $codreg
$ws.Cells.Item($r, 6).Value2 = $tipo
$ws.Cells.Item($r, 7).Value2 = $productoId
$ws.Cells.Item($r, 8).Value2 = $producto
$ws.Cells.Item($r, 9).Value2 = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = $variedad
$ws.Cells.Item($r, 12).Value2 = "Primera"
$ws.Cells.Item($r, 13).Value2 = 300
$ws.Cells.Item($r, 14).Value2 = 17000
$ws.Cells.Item($r, 15).Value2 = 17000
$ws.Cells.Item($r, 16).Value2 = 17000
$ws.Cells.Item($r, 17).Value2 = $unidad
$ws.Cells.Item($r, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value2 = 1133
$ws.Cells.Item($r, 20).Value2 = $kgUnidad

# Row 131: Segunda
$r = 131
$ws.Cells.Item($r, 1).Value2 = $mercadoId
$ws.Cells.Item($r, 2).Value2 = $mercado
$ws.Cells.Item($r, 3).Value2 = $region
$ws.Cells.Item($r, 4).Value2 = $fecha
$ws.Cells.Item($r, 5).Value2 = $codreg
$ws.Cells.Item($r, 6).Value2 = $tipo
$ws.Cells.Item($r, 7).Value2 = $productoId
$ws.Cells.Item($r, 8).Value2 = $producto
$ws.Cells.Item($r, 9).Value2 = $categoriaId
$ws.Cells.Item($r, 10).Value2 = $categoria
$ws.Cells.Item($r, 11).Value2 = $variedad
$ws.Cells.Item($r, 12).Value2 = "Segunda"
$ws.Cells.Item($r, 13).Value2 = 300
$ws.Cells.Item($r, 14).Value2 = 15000
$ws.Cells.Item($r, 15).Value2 = 15000
$ws.Cells.Item($r, 16).Value2 = 15000
$ws.Cells.Item($r, 17).Value2 = $unidad
$ws.Cells.Item($r, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value2 = 1000
$ws.Cells.Item($r, 20).Value2 = $kgUnidad
